$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new column before column B ("Gender"). This shifts the former
#    B:L columns (Firstname..ReponseSecrete) one place to the right (C:M),
#    carrying cell styles and column widths along with them.
# ---------------------------------------------------------------------------
$ws.Columns("B:B").Insert()

# ---------------------------------------------------------------------------
# 2) Populate the new "Gender" column (header + 8 data rows).
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "Gender"

$genders = @("F","M","M","M","M","F","M","M")
for ($i = 0; $i -lt $genders.Length; $i++) {
    $row = $i + 2
    $ws.Range("B" + $row).Value = $genders[$i]
}

# ---------------------------------------------------------------------------
# 3) The hyperlinks that used to live in column D (Email) are now in column
#    E, but this runtime does not automatically re-anchor Hyperlink ranges
#    when a column is inserted, so rebuild them explicitly: capture the
#    e-mail addresses shown in the cells, drop every existing hyperlink,
#    then re-add them anchored on the shifted column E.
# ---------------------------------------------------------------------------
$emails = @()
for ($row = 2; $row -le 9; $row++) {
    $emails += $ws.Range("E" + $row).Value()
}

$ws.Hyperlinks.Delete()

for ($i = 0; $i -lt $emails.Length; $i++) {
    $row = $i + 2
    $target = $ws.Range("E" + $row)
    $ws.Hyperlinks.Add($target, "mailto:" + $emails[$i], "", "", $emails[$i])
}

# ---------------------------------------------------------------------------
# 4) Append the new "IsInfo" and "Lang" columns (N and O) at the end.
# ---------------------------------------------------------------------------
$ws.Range("N1").Value = "IsInfo"
$ws.Range("O1").Value = "Lang"

$isInfo = @(1, 1, 0, 1, 1, 1, 0, 1)
$lang   = @("fr", "fr", "fr", "fr", "fr", "fr", "en", "en")
for ($i = 0; $i -lt $isInfo.Length; $i++) {
    $row = $i + 2
    $ws.Range("N" + $row).Value = $isInfo[$i]
    $ws.Range("O" + $row).Value = $lang[$i]
}

# ---------------------------------------------------------------------------
# 5) Restore the active selection shown in the saved workbook.
# ---------------------------------------------------------------------------
$ws.Range("B9").Select()
